$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header D1 from "PpMeHg" to "pMeHg" (column E "RpMeHg" will be removed below)
$ws.Range("D1").Value = "pMeHg"

# Update the D4/D5 pMeHg values (previously flagged as "<0.1" via columns D/E, now simplified to 0.05)
$ws.Range("D4").Value = 0.05
$ws.Range("D5").Value = 0.05

# Remove the now-unused RpMeHg column (E) entirely, which also removes the "<" flag cells E4/E5
$ws.Columns("E").Delete()

# D15 now uses the repurposed 2-decimal numeric style
$ws.Range("D15").NumberFormat = "0.00"

# Update the active selection to reflect the new last-used column (E31, one past column D)
[void]$ws.Range("E31").Select()
